$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab (date moved from 12-26 to 12-27)
$ws.Name = "Through 2021-12-27"

# Row 10 (August) - 2021 column updates
$ws.Range("T10").Value = 9
$ws.Range("U10").Value = 151
$ws.Range("V10").Value = 0.0562

# Row 14 (December) label update
$ws.Range("A14").Value = "December (through 12-27)"

# Row 14 (December) - 2015 columns
$ws.Range("C14").Value = 36
$ws.Range("D14").Value = 0.1

# Row 14 (December) - 2016 columns
$ws.Range("F14").Value = 79
$ws.Range("G14").Value = 0.092

# Row 14 (December) - 2017 columns
$ws.Range("I14").Value = 94
$ws.Range("J14").Value = 0.1132

# Row 14 (December) - 2018 columns
$ws.Range("L14").Value = 59
$ws.Range("M14").Value = 0.0781

# Row 14 (December) - 2020 columns
$ws.Range("R14").Value = 118
$ws.Range("S14").Value = 0.0635

# Row 14 (December) - 2021 columns
$ws.Range("U14").Value = 164
$ws.Range("V14").Value = 0.012

# Row 15 (Total) - 2015 columns
$ws.Range("C15").Value = 294
$ws.Range("D15").Value = 0.1118

# Row 15 (Total) - 2016 columns
$ws.Range("F15").Value = 583
$ws.Range("G15").Value = 0.1031

# Row 15 (Total) - 2017 columns
$ws.Range("I15").Value = 852
$ws.Range("J15").Value = 0.0809

# Row 15 (Total) - 2018 columns
$ws.Range("L15").Value = 667
$ws.Range("M15").Value = 0.1059

# Row 15 (Total) - 2020 columns
$ws.Range("R15").Value = 1318
$ws.Range("S15").Value = 0.0518

# Row 15 (Total) - 2021 columns
$ws.Range("T15").Value = 103
$ws.Range("U15").Value = 1706
$ws.Range("V15").Value = 0.0569
